$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Range("G30").Value = [double]"-1.14473713164009"
$ws.Range("H30").Value = [double]"0.26055311655631"

# Row 31
$ws.Range("C31").Value = "manual_few_shot"
$ws.Range("F31").Value = [double]"0.3487301702745721"
$ws.Range("G31").Value = [double]"-0.9872890808083997"
$ws.Range("H31").Value = [double]"0.3292918121132551"

# Row 32
$ws.Range("C32").Value = "self_refine"
$ws.Range("F32").Value = [double]"0.351156728855353"
$ws.Range("G32").Value = [double]"-1.829021014606149"
$ws.Range("H32").Value = [double]"0.07672940280202731"

# Row 33
$ws.Range("C33").Value = "tree_of_thought"
$ws.Range("F33").Value = [double]"0.3478996421926851"
$ws.Range("G33").Value = [double]"-1.111119966153743"
$ws.Range("H33").Value = [double]"0.2829354362531543"

# Row 34
$ws.Range("C34").Value = "zero_shot_cot"
$ws.Range("F34").Value = [double]"0.3441867682123834"
$ws.Range("G34").Value = [double]"-1.298666300051236"
$ws.Range("H34").Value = [double]"0.2018843348096545"

# Row 35
$ws.Range("C35").Value = "least_to_most"
$ws.Range("F35").Value = [double]"0.3498959518962656"
$ws.Range("G35").Value = [double]"-2.408247956271698"
$ws.Range("H35").Value = [double]"0.02196460748337448"
$ws.Range("I35").Value = "Yes"

# Row 36
$ws.Range("C36").Value = "manual_cot"
$ws.Range("F36").Value = [double]"0.3551923713253787"
$ws.Range("G36").Value = [double]"-2.374394674574453"
$ws.Range("H36").Value = [double]"0.02374316969700669"
$ws.Range("I36").Value = "Yes"

# Row 44
$ws.Range("C44").Value = "least_to_most"
$ws.Range("F44").Value = [double]"0.3688316461987508"
$ws.Range("G44").Value = [double]"-0.4098421582418101"
$ws.Range("H44").Value = [double]"0.6877150738021158"

# Row 45
$ws.Range("C45").Value = "manual_cot"
$ws.Range("F45").Value = [double]"0.3506891984872582"
$ws.Range("G45").Value = [double]"0.2506884562793319"
$ws.Range("H45").Value = [double]"0.8076843000403842"

# Row 46
$ws.Range("C46").Value = "manual_few_shot"
$ws.Range("F46").Value = [double]"0.363439672985659"
$ws.Range("G46").Value = [double]"-1.099920436908166"
$ws.Range("H46").Value = [double]"0.2814470372062072"

# Row 47
$ws.Range("C47").Value = "self_refine"
$ws.Range("F47").Value = [double]"0.3927645311833518"
$ws.Range("G47").Value = [double]"-1.881009639850796"
$ws.Range("H47").Value = [double]"0.07624865825236271"

# Row 48
$ws.Range("C48").Value = "tree_of_thought"
$ws.Range("F48").Value = [double]"0.3970522049201494"
$ws.Range("G48").Value = [double]"-0.178304072906799"
$ws.Range("H48").Value = [double]"0.8876686001034875"

# Row 49
$ws.Range("C49").Value = "zero_shot_cot"
$ws.Range("F49").Value = [double]"0.3859400964114446"
$ws.Range("G49").Value = [double]"0.09440419099120144"
$ws.Range("H49").Value = [double]"0.9255118747930475"

# Row 50
$ws.Range("C50").Value = "ape_zero_shot_cot"
$ws.Range("F50").Value = [double]"0.3910591202955087"
$ws.Range("G50").Value = [double]"-2.106654550723698"
$ws.Range("H50").Value = [double]"0.04734782328834347"
$ws.Range("I50").Value = "Yes"

# Row 61
$ws.Range("G61").Value = [double]"0.380597446300002"
$ws.Range("H61").Value = [double]"0.7049950980741344"

# Row 72
$ws.Range("C72").Value = "self_refine"
$ws.Range("F72").Value = [double]"0.3655432912693316"
$ws.Range("G72").Value = [double]"-0.1770794927771922"
$ws.Range("H72").Value = [double]"0.859867225764093"

# Row 73
$ws.Range("C73").Value = "tree_of_thought"
$ws.Range("F73").Value = [double]"0.4301229604135541"
$ws.Range("G73").Value = [double]"-1.462963211730259"
$ws.Range("H73").Value = [double]"0.1550220395534408"

# Row 74
$ws.Range("C74").Value = "ape_zero_shot_cot"
$ws.Range("F74").Value = [double]"0.4195166577919803"
$ws.Range("G74").Value = [double]"-2.722508804502758"
$ws.Range("H74").Value = [double]"0.007912598988651227"
$ws.Range("I74").Value = "Yes"

# Row 75
$ws.Range("C75").Value = "least_to_most"
$ws.Range("F75").Value = [double]"0.4795389905014912"
$ws.Range("G75").Value = [double]"-5.704818886801226"
$ws.Range("H75").Value = [double]"1.37328393154626e-07"
$ws.Range("I75").Value = "Yes"

# Row 76
$ws.Range("C76").Value = "manual_cot"
$ws.Range("F76").Value = [double]"0.4786850615266719"
$ws.Range("G76").Value = [double]"-5.848788038540913"
$ws.Range("H76").Value = [double]"7.151516652751101e-08"
$ws.Range("I76").Value = "Yes"

# Row 77
$ws.Range("C77").Value = "manual_few_shot"
$ws.Range("F77").Value = [double]"0.4762447813337619"
$ws.Range("G77").Value = [double]"-5.402586153841306"
$ws.Range("H77").Value = [double]"4.852185926605229e-07"
$ws.Range("I77").Value = "Yes"

# Row 78
$ws.Range("G78").Value = [double]"-2.345575437840131"
$ws.Range("H78").Value = [double]"0.02150696795521951"
$ws.Range("I78").Value = "Yes"

# Row 86
$ws.Range("G86").Value = [double]"-0.5706359002568039"
$ws.Range("H86").Value = [double]"0.5721119861438464"

# Row 87
$ws.Range("G87").Value = [double]"-0.04573220039524427"
$ws.Range("H87").Value = [double]"0.9638079071817984"

# Row 88
$ws.Range("G88").Value = [double]"-0.5956350360123538"
$ws.Range("H88").Value = [double]"0.5556087468580261"

# Row 89
$ws.Range("C89").Value = "self_refine"
$ws.Range("F89").Value = [double]"0.3976991678288766"
$ws.Range("G89").Value = [double]"0.1304041786778109"
$ws.Range("H89").Value = [double]"0.8970629935145363"

# Row 90
$ws.Range("C90").Value = "tree_of_thought"
$ws.Range("F90").Value = [double]"0.4552911088706001"
$ws.Range("G90").Value = [double]"1.232568738202675"
$ws.Range("H90").Value = [double]"0.2355418166853398"

# Row 91
$ws.Range("C91").Value = "manual_few_shot"
$ws.Range("F91").Value = [double]"0.3657057845474713"
$ws.Range("G91").Value = [double]"2.658636260743372"
$ws.Range("H91").Value = [double]"0.01113814179210349"
$ws.Range("I91").Value = "Yes"

# Row 92
$ws.Range("G92").Value = [double]"-2.51052457717059"
$ws.Range("H92").Value = [double]"0.01643009978264356"
$ws.Range("I92").Value = "Yes"

# Row 100
$ws.Range("G100").Value = [double]"-0.5398998232920198"
$ws.Range("H100").Value = [double]"0.5949457001656505"

# Row 101
$ws.Range("G101").Value = [double]"-0.09070866577678419"
$ws.Range("H101").Value = [double]"0.9289243014000621"

# Row 102
$ws.Range("G102").Value = [double]"-0.68775779570202"
$ws.Range("H102").Value = [double]"0.5089404360165983"

# Row 103
$ws.Range("C103").Value = "self_refine"
$ws.Range("F103").Value = [double]"0.3403139039125427"
$ws.Range("G103").Value = [double]"0.6324113498329108"
$ws.Range("H103").Value = [double]"0.5350696318477244"

# Row 104
$ws.Range("C104").Value = "tree_of_thought"
$ws.Range("F104").Value = [double]"0.223879186997216"
$ws.Range("G104").Value = [double]"1.241549692571777"
$ws.Range("H104").Value = [double]"0.4316614758653485"

# Row 105
$ws.Range("C105").Value = "zero_shot_cot"
$ws.Range("F105").Value = [double]"0.3573737035024921"
$ws.Range("G105").Value = [double]"1.439058993114993"
$ws.Range("H105").Value = [double]"0.1620648744709412"

# Row 106
$ws.Range("C106").Value = "manual_few_shot"
$ws.Range("F106").Value = [double]"0.433228101589585"
$ws.Range("G106").Value = [double]"-2.928667469841096"
$ws.Range("H106").Value = [double]"0.006993274462604279"
$ws.Range("I106").Value = "Yes"

# Row 226
$ws.Range("G226").Value = [double]"-29.84962311319859"
$ws.Range("H226").Value = [double]"2.597549811426689e-10"
$ws.Range("I226").Value = "Yes"

# Row 227
$ws.Range("G227").Value = [double]"-26.70571328295164"
$ws.Range("H227").Value = [double]"7.000859820572747e-10"
$ws.Range("I227").Value = "Yes"

# Row 228
$ws.Range("G228").Value = [double]"-2.954195783503986"
$ws.Range("H228").Value = [double]"0.01611071653191161"
$ws.Range("I228").Value = "Yes"

# Row 229
$ws.Range("G229").Value = [double]"-11.22497216032182"
$ws.Range("H229").Value = [double]"1.356959540831356e-06"
$ws.Range("I229").Value = "Yes"

# Row 230
$ws.Range("G230").Value = [double]"-9.221981556055331"
$ws.Range("H230").Value = [double]"6.992072395885232e-06"
$ws.Range("I230").Value = "Yes"

# Row 231
$ws.Range("G231").Value = [double]"-11.75894243853278"
$ws.Range("H231").Value = [double]"9.151111215642479e-07"
$ws.Range("I231").Value = "Yes"

# Row 232
$ws.Range("G232").Value = [double]"-40.41658075592244"
$ws.Range("H232").Value = [double]"1.73010492641202e-11"
$ws.Range("I232").Value = "Yes"
